$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("I7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 107"
